$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date in column C for rows 2-12
# from 45208 (2023-10-09) to 45212 (2023-10-13)
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = 45212
}
